$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 87.95
$ws.Range("I2").Value = 47.266666
$ws.Range("J2").Value = 210
$ws.Range("K2").Value = 47.266666
$ws.Range("L2").Value = 210
$ws.Range("M2").Value = 65.733334
$ws.Range("N2").Value = -436

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1547006.5
$ws.Range("I32").Value = 1585537.8
$ws.Range("K32").Value = 1585537.8
$ws.Range("M32").Value = -1585250.8
$ws.Range("H61").Value = 27785676
$ws.Range("I61").Value = 5293.2085
$ws.Range("J61").Value = 83346440
$ws.Range("K61").Value = 5293.2085
$ws.Range("L61").Value = 83346440
$ws.Range("M61").Value = -5081.2085
$ws.Range("N61").Value = -83346864
$ws.Range("H74").Value = 29453.916
$ws.Range("I74").Value = 36073.062
$ws.Range("J74").Value = 5523.154
$ws.Range("K74").Value = 36073.062
$ws.Range("L74").Value = 5523.154
$ws.Range("M74").Value = -35199.062
$ws.Range("N74").Value = -7271.154
$ws.Range("H77").Value = 29453.916
$ws.Range("I77").Value = 36073.062
$ws.Range("J77").Value = 5523.154
$ws.Range("K77").Value = 180365.31
$ws.Range("L77").Value = 27615.77
$ws.Range("M77").Value = -175997.31
$ws.Range("N77").Value = -36351.77
$ws.Range("H132").Value = 5224.1836
$ws.Range("I132").Value = 2630.724
$ws.Range("K132").Value = 7892.172
$ws.Range("M132").Value = -5362.172
$ws.Range("H136").Value = 27785676
$ws.Range("I136").Value = 5293.2085
$ws.Range("J136").Value = 83346440
$ws.Range("K136").Value = 15879.6255
$ws.Range("L136").Value = 250039320
$ws.Range("M136").Value = -13329.6255
$ws.Range("N136").Value = -250044420

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1594.4242
$ws.Range("I94").Value = 1042.4286
$ws.Range("J94").Value = 2560.4167
$ws.Range("K94").Value = 1042.4286
$ws.Range("L94").Value = 2560.4167
$ws.Range("M94").Value = -591.4286
$ws.Range("N94").Value = -3462.4167
$ws.Range("H134").Value = 6255855
$ws.Range("I134").Value = 11906210
$ws.Range("J134").Value = 10726.053
$ws.Range("K134").Value = 35718630
$ws.Range("L134").Value = 32178.159
$ws.Range("M134").Value = -35716095
$ws.Range("N134").Value = -37248.159

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5274.4287
$ws.Range("I16").Value = 966.4
$ws.Range("K16").Value = 966.4
$ws.Range("M16").Value = -679.4
$ws.Range("H31").Value = 7414820
$ws.Range("I31").Value = 2437.9167
$ws.Range("J31").Value = 10110231
$ws.Range("K31").Value = 2437.9167
$ws.Range("L31").Value = 10110231
$ws.Range("M31").Value = -2142.9167
$ws.Range("N31").Value = -10110821
$ws.Range("H34").Value = 7414820
$ws.Range("I34").Value = 2437.9167
$ws.Range("J34").Value = 10110231
$ws.Range("K34").Value = 2437.9167
$ws.Range("L34").Value = 10110231
$ws.Range("M34").Value = -2235.9167
$ws.Range("N34").Value = -10110635
$ws.Range("H58").Value = 5624.44
$ws.Range("I58").Value = 2002.3043
$ws.Range("J58").Value = 8709.963
$ws.Range("K58").Value = 2002.3043
$ws.Range("L58").Value = 8709.963
$ws.Range("M58").Value = -1799.3043
$ws.Range("N58").Value = -9115.963
$ws.Range("H99").Value = 7527.846
$ws.Range("I99").Value = 7384.4614
$ws.Range("J99").Value = 7671.231
$ws.Range("K99").Value = 7384.4614
$ws.Range("L99").Value = 7671.231
$ws.Range("M99").Value = -5886.4614
$ws.Range("N99").Value = -10667.231
$ws.Range("H113").Value = 5274.4287
$ws.Range("I113").Value = 966.4
$ws.Range("K113").Value = 966.4
$ws.Range("M113").Value = 1203.6
$ws.Range("H122").Value = 3048.318
$ws.Range("I122").Value = 1957.5714
$ws.Range("K122").Value = 5872.7142
$ws.Range("M122").Value = -3422.7142
$ws.Range("H126").Value = 7527.846
$ws.Range("I126").Value = 7384.4614
$ws.Range("J126").Value = 7671.231
$ws.Range("K126").Value = 22153.3842
$ws.Range("L126").Value = 23013.693
$ws.Range("M126").Value = -19683.3842
$ws.Range("N126").Value = -27953.693
$ws.Range("H132").Value = 6156941.5
$ws.Range("I132").Value = 1085.6875
$ws.Range("J132").Value = 23538182
$ws.Range("K132").Value = 3257.0625
$ws.Range("L132").Value = 70614546
$ws.Range("M132").Value = -727.0625
$ws.Range("N132").Value = -70619606
$ws.Range("H136").Value = 5624.44
$ws.Range("I136").Value = 2002.3043
$ws.Range("J136").Value = 8709.963
$ws.Range("K136").Value = 6006.9129
$ws.Range("L136").Value = 26129.889
$ws.Range("M136").Value = -3456.9129
$ws.Range("N136").Value = -31229.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 6797.885
$ws.Range("I132").Value = 2316.8
$ws.Range("K132").Value = 20851.2
$ws.Range("M132").Value = -18321.2
$ws.Range("H134").Value = 56162.2
$ws.Range("I134").Value = 88452
$ws.Range("J134").Value = 7727.5
$ws.Range("K134").Value = 265356
$ws.Range("L134").Value = 23182.5
$ws.Range("M134").Value = -260286
$ws.Range("N134").Value = -33322.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6378.6665
$ws.Range("I70").Value = 4101.067
$ws.Range("J70").Value = 8276.666999999999
$ws.Range("K70").Value = 4101.067
$ws.Range("L70").Value = 8276.666999999999
$ws.Range("M70").Value = -3831.067
$ws.Range("N70").Value = -8816.666999999999
$ws.Range("H73").Value = 6378.6665
$ws.Range("I73").Value = 4101.067
$ws.Range("J73").Value = 8276.666999999999
$ws.Range("K73").Value = 4101.067
$ws.Range("L73").Value = 8276.666999999999
$ws.Range("M73").Value = -3165.067
$ws.Range("N73").Value = -10148.667
$ws.Range("H102").Value = 2363.7896
$ws.Range("I102").Value = 2463.8572
$ws.Range("J102").Value = 1750.875
$ws.Range("K102").Value = 2463.8572
$ws.Range("L102").Value = 1750.875
$ws.Range("M102").Value = -841.8571999999999
$ws.Range("N102").Value = -4994.875
$ws.Range("H132").Value = 3400.907
$ws.Range("I132").Value = 1228.5454
$ws.Range("J132").Value = 10569.7
$ws.Range("K132").Value = 3685.6362
$ws.Range("L132").Value = 31709.1
$ws.Range("M132").Value = -1155.6362
$ws.Range("N132").Value = -36769.10000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1198.5
$ws.Range("I16").Value = 1198.5
$ws.Range("K16").Value = 1198.5
$ws.Range("M16").Value = -1028.5
$ws.Range("H132").Value = 6333628
$ws.Range("I132").Value = 10871688
$ws.Range("K132").Value = 32615064
$ws.Range("M132").Value = -32612534
$ws.Range("H136").Value = 6171.0894
$ws.Range("I136").Value = 2149.1
$ws.Range("K136").Value = 6447.299999999999
$ws.Range("M136").Value = -3897.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15628378
$ws.Range("I132").Value = 18870996
$ws.Range("J132").Value = 4858.909
$ws.Range("K132").Value = 56612988
$ws.Range("L132").Value = 14576.727
$ws.Range("M132").Value = -56610458
$ws.Range("N132").Value = -19636.727
$ws.Range("H136").Value = 27301722
$ws.Range("I136").Value = 45456040
$ws.Range("J136").Value = 675390.7
$ws.Range("K136").Value = 136368120
$ws.Range("L136").Value = 2026172.1
$ws.Range("M136").Value = -136365570
$ws.Range("N136").Value = -2031272.1
